# chore: adapt column header formatting to respective input file names (#7)
#
# Renames the "_old"/"_new" column-header suffixes used for the AHB diff
# columns to the concrete format-version names ("_FV2410" / "_FV2504"),
# wraps the sheet's data range in a native Excel Table (ListObject) and
# freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells (A1:J1 -> *_FV2410, L1:U1 -> *_FV2504) -----
$fv2410 = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
$fv2504 = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $fv2410.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2410[$i]
}
# Column 11 ("K") holds the "diff" header and is left untouched.
for ($i = 0; $i -lt $fv2504.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2504[$i]
}

# --- 2. Turn the used range into a native Excel table ----------------------
$dataRange = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ----------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Output "edit.ps1: headers renamed, Table1 added, header row frozen"
